# This script swaps the contents of columns D ("codeforiati:group-name")
# and E ("codeforiati:group-code") for every row of the worksheet,
# including the header row. This mirrors the upstream codeforIATI
# codelists rebuild, which reordered the "group-code" / "group-name"
# columns so that group-code now precedes group-name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work out how many rows are actually populated.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) {
    $lastRow = 1
}

$dRange = $ws.Range("D1:D$lastRow")
$eRange = $ws.Range("E1:E$lastRow")

# Value2 returns the raw values (no currency/date formatting) as a
# 2D array when the range spans multiple cells.
$dValues = $dRange.Value2
$eValues = $eRange.Value2

# Swap the two columns' contents in one shot.
$dRange.Value = $eValues
$eRange.Value = $dValues
